$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextCell 'D2' '59.064.25'
Set-TextCell 'E2' '  +0.94%  '

Set-TextCell 'D3' '2.630.13'
Set-TextCell 'E3' '  +3.45%  '

Set-TextCell 'E4' '  +0.24%  '

Set-TextCell 'D5' '519.00'
Set-TextCell 'E5' '  +2.46%  '

Set-TextCell 'D6' '145.65'
Set-TextCell 'E6' '  +1.83%  '

Set-TextCell 'E7' '  -0.03%  '

Set-TextCell 'E8' '  +0.31%  '

Set-TextCell 'D9' '2.656.08'
Set-TextCell 'E9' '  +4.18%  '

Set-TextCell 'D10' '6.28'
Set-TextCell 'E10' '  +3.19%  '

Set-TextCell 'E11' '  +3.16%  '

Set-TextCell 'E12' '  +1.66%  '

Set-TextCell 'E13' '  -1.61%  '

Set-TextCell 'D14' '3.120.74'
Set-TextCell 'E14' '  +4.50%  '

Set-TextCell 'D15' '59.084.36'
Set-TextCell 'E15' '  +1.02%  '

Set-TextCell 'E16' '  +1.56%  '

Set-TextCell 'E17' '  +1.58%  '

Set-TextCell 'D18' '2.658.24'
Set-TextCell 'E18' '  +4.50%  '

Set-TextCell 'D19' '350.26'
Set-TextCell 'E19' '  +3.26%  '

Set-TextCell 'E20' '  -0.15%  '

Set-TextCell 'D21' '10.34'
Set-TextCell 'E21' '  +2.59%  '

Set-TextCell 'E22' '  +3.83%  '

Set-TextCell 'D23' '0.998'
Set-TextCell 'E23' '  -0.20%  '

Set-TextCell 'D24' '62.10'
Set-TextCell 'E24' '  +2.37%  '

Set-TextCell 'E25' '  +1.93%  '

Set-TextCell 'D26' '2.761.49'
Set-TextCell 'E26' '  +4.16%  '

Set-TextCell 'E27' '  +2.73%  '

Set-TextCell 'D28' '0.998'
Set-TextCell 'E28' '  -0.07%  '

Set-TextCell 'D29' '0.0₃0805'
Set-TextCell 'E29' '  +2.41%  '

Set-TextCell 'E30' '  +2.54%  '

Set-TextCell 'D31' '0.999'
Set-TextCell 'E31' '  +0.00%  '

Set-TextCell 'D32' '6.25'
Set-TextCell 'E32' '  +7.20%  '

Set-TextCell 'D33' '19.00'
Set-TextCell 'E33' '  +2.67%  '

Set-TextCell 'E34' '  +2.99%  '

Set-TextCell 'D35' '149.93'
Set-TextCell 'E35' '  +0.19%  '

Set-TextCell 'D36' '0.961'
Set-TextCell 'E36' '  +5.38%  '

Set-TextCell 'E37' '  +3.42%  '

Set-TextCell 'E38' '  +2.43%  '

Set-TextCell 'D39' '36.75'
Set-TextCell 'E39' '  +1.87%  '

Set-TextCell 'D40' '0.844'
Set-TextCell 'E40' '  +2.82%  '

Set-TextCell 'D41' '3.69'
Set-TextCell 'E41' '  +4.88%  '

Set-TextCell 'E42' '  +2.10%  '

Set-TextCell 'D43' '277.21'
Set-TextCell 'E43' '  -1.83%  '

Set-TextCell 'D44' '0.995'
Set-TextCell 'E44' '  -0.23%  '

Set-TextCell 'B45' 'Mantle'
Set-TextCell 'C45' 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextCell 'D45' '0.610'
Set-TextCell 'E45' '  +1.67%  '

Set-TextCell 'B46' 'Stellar'
Set-TextCell 'C46' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell 'D46' '0.0986'
Set-TextCell 'E46' '  -1.15%  '

Set-TextCell 'D47' '19.64'
Set-TextCell 'E47' '  +5.54%  '

Set-TextCell 'D48' '0.0521'
Set-TextCell 'E48' '  -2.01%  '

Set-TextCell 'B49' 'WhiteBITCoin'
Set-TextCell 'C49' 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextCell 'D49' '10.30'
Set-TextCell 'E49' '  +0.13%  '

Set-TextCell 'B50' 'VeChain'
Set-TextCell 'C50' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell 'D50' '0.0230'
Set-TextCell 'E50' '  +1.53%  '

Set-TextCell 'B51' 'Maker'
Set-TextCell 'C51' 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextCell 'D51' '1.986.79'
Set-TextCell 'E51' '  +4.85%  '
